$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column I that need a value of 0 added. Rows 21, 23, 25, 30, 34 and
# 36 already carry a real "FRA discount" figure in column I and must be
# left untouched; every other data row (2-37) gets a 0.
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,24,26,27,28,29,31,32,33,35,37)

# Column H (entry fee) already carries this currency number format; reusing
# the exact format-code string (quotes/backslash included) keeps the newly
# written cells on that same existing style record instead of minting a
# duplicate, near-identical one.
$currencyFormat = '"£"#,##0;[Red]\-"£"#,##0'

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.NumberFormat = $currencyFormat
    $cell.Value = 0
}

# Scroll the frozen view down so row 10 is the first visible row below the
# frozen header, and leave the selection on I38 (one row below the last
# data row), matching the saved view state of the edited workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$ws.Range("I38").Select()
